# Update the "PSF Shifts" worksheet (sheet3) with new FDPR run results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PSF Shifts")
$ws.Activate()

# --- Row 2 ---
$ws.Range("A2").Value = 45058.67454861111
$ws.Range("B2").Value = 45058.84121527777
$ws.Range("D2").Value = -61
$ws.Range("E2").Value = -62
$ws.Range("F2").Value = 228.9
$ws.Range("G2").Value = 321.3
$ws.Range("H2").Value = 1.3
$ws.Range("I2").Value = 27.1
$ws.Range("J2").Value = 13
$ws.Range("K2").Value = 271
$ws.Range("L2").Value = -0.0045
$ws.Range("M2").Value = -4.9981
$ws.Range("N2").Value = 0.8983
$ws.Range("O2").Value = 0.6421
$ws.Range("P2").Value = 0.2848
$ws.Range("Q2").Value = 0.2853
$ws.Range("R2").Value = 0.0164
$ws.Range("S2").Value = 0.0492
$ws.Range("T2").Value = "PR1-830-20230512_161121.A1.sol.h5"

# --- Row 3 ---
$ws.Range("A3").Value = 45058.68513888889
$ws.Range("B3").Value = 45058.85180555555
$ws.Range("D3").Value = -32
$ws.Range("E3").Value = -32
$ws.Range("F3").Value = 214.2
$ws.Range("G3").Value = 318.9
$ws.Range("H3").Value = -1.1
$ws.Range("I3").Value = 41.8
$ws.Range("J3").Value = -11
$ws.Range("K3").Value = 418
$ws.Range("L3").Value = -0.0027
$ws.Range("M3").Value = -4.998
$ws.Range("N3").Value = 111.1702
$ws.Range("O3").Value = 23.1252
$ws.Range("P3").Value = -46.2409
$ws.Range("Q3").Value = -46.2407
$ws.Range("R3").Value = -4.886
$ws.Range("S3").Value = -10.5802
$ws.Range("T3").Value = "PR3-830-20230512_162636.A1.sol.h5"

# --- Row 4 ---
$ws.Range("A4").Value = 45058.69462962963
$ws.Range("B4").Value = 45058.86129629629
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 21
$ws.Range("F4").Value = 208.8
$ws.Range("G4").Value = 322
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 47.2
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 472
$ws.Range("L4").Value = -0.0008
$ws.Range("M4").Value = -5.0028
$ws.Range("N4").Value = 102.96
$ws.Range("O4").Value = 31.7847
$ws.Range("P4").Value = 58.8434
$ws.Range("Q4").Value = 58.8439
$ws.Range("R4").Value = 5.2803
$ws.Range("S4").Value = -10.3049
$ws.Range("T4").Value = "PR4-830-20230512_164016.A1.sol.h5"

# --- Row 5 ---
$ws.Range("A5").Value = 45058.70511574074
$ws.Range("B5").Value = 45058.8717824074
$ws.Range("D5").Value = -174
$ws.Range("E5").Value = -203
$ws.Range("F5").Value = 238.1
$ws.Range("G5").Value = 320.9
$ws.Range("H5").Value = 0.9
$ws.Range("I5").Value = 17.9
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 179
$ws.Range("L5").Value = 0.0008
$ws.Range("M5").Value = -4.9995
$ws.Range("N5").Value = -101.4177
$ws.Range("O5").Value = 31.8065
$ws.Range("P5").Value = 58.4249
$ws.Range("Q5").Value = 58.4265
$ws.Range("R5").Value = 5.2705
$ws.Range("S5").Value = 10.4165
$ws.Range("T5").Value = "PR5-830-20230512_165522.A1.sol.h5"

# --- Row 6 ---
$ws.Range("A6").Value = 45058.71503472222
$ws.Range("B6").Value = 45058.88170138889
$ws.Range("D6").Value = -103
$ws.Range("E6").Value = -106
$ws.Range("F6").Value = 226
$ws.Range("G6").Value = 322.5
$ws.Range("H6").Value = 2.5
$ws.Range("I6").Value = 30
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 0.0022
$ws.Range("M6").Value = -5.0018
$ws.Range("N6").Value = -109.4993
$ws.Range("O6").Value = 23.3409
$ws.Range("P6").Value = -46.7926
$ws.Range("Q6").Value = -46.7934
$ws.Range("R6").Value = -4.9126
$ws.Range("S6").Value = 10.7008
$ws.Range("T6").Value = "PR2-830-20230512_170939.A1.sol.h5"

# The "PR results file" header (T1) switches from centered to left-aligned text,
# which introduces a new cell style in the workbook.
$ws.Range("T1").HorizontalAlignment = -4131

# The trailing blank row (row 11) is removed, shrinking the sheet's used range
# from A1:T11 down to A1:T10.
$ws.Rows("11:11").Delete()

# Update the selected cell shown when the sheet is reopened.
$ws.Range("D5").Select()
